$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.817.58'
$ws.Range('E2').Value = '  -9.88%  '
$ws.Range('D3').Value = '2.361.59'
$ws.Range('E3').Value = '  -12.85%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '454.57'
$ws.Range('E5').Value = '  -9.87%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '128.77'
$ws.Range('E6').Value = '  -7.90%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.476'
$ws.Range('E8').Value = '  -9.96%  '
$ws.Range('D9').Value = '2.378.67'
$ws.Range('E9').Value = '  -12.39%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0935'
$ws.Range('E10').Value = '  -9.28%  '
$ws.Range('E11').Value = '  -14.34%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.309'
$ws.Range('E12').Value = '  -10.66%  '
$ws.Range('E13').Value = '  -4.77%  '
$ws.Range('D14').Value = '2.778.52'
$ws.Range('E14').Value = '  -12.87%  '
$ws.Range('D15').Value = '52.956.14'
$ws.Range('E15').Value = '  -9.48%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '19.37'
$ws.Range('E16').Value = '  -9.61%  '
$ws.Range('E17').Value = '  -4.48%  '
$ws.Range('D18').Value = '2.380.82'
$ws.Range('E18').Value = '  -12.14%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.12'
$ws.Range('E19').Value = '  -12.29%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '303.13'
$ws.Range('E20').Value = '  -11.15%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.25'
$ws.Range('E21').Value = '  -14.89%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.992'
$ws.Range('E22').Value = '  -0.78%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.66'
$ws.Range('E23').Value = '  +1.02%  '
$ws.Range('E24').Value = '  -14.94%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '55.24'
$ws.Range('E25').Value = '  -11.70%  '
$ws.Range('E26').Value = '  +1.10%  '
$ws.Range('E27').Value = '  -10.35%  '
$ws.Range('D28').Value = '2.486.98'
$ws.Range('E28').Value = '  -11.91%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.149'
$ws.Range('E29').Value = '  -12.30%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.97'
$ws.Range('E30').Value = '  -5.40%  '
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('D32').Value = '0.0₃0706'
$ws.Range('E32').Value = '  -14.56%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '144.68'
$ws.Range('E33').Value = '  -2.39%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '17.40'
$ws.Range('E34').Value = '  -8.30%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.40'
$ws.Range('E35').Value = '  -12.49%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.90'
$ws.Range('E36').Value = '  -7.60%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.45'
$ws.Range('E37').Value = '  -16.82%  '
$ws.Range('E38').Value = '  -7.37%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.779'
$ws.Range('E39').Value = '  -16.32%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.997'
$ws.Range('E40').Value = '  +0.18%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '32.79'
$ws.Range('E41').Value = '  -8.76%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.584'
$ws.Range('E42').Value = '  -2.01%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.22'
$ws.Range('E43').Value = '  -7.67%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0517'
$ws.Range('E44').Value = '  -5.89%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '10.14'
$ws.Range('E45').Value = '  -1.96%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.21'
$ws.Range('E46').Value = '  -11.62%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '1.921.19'
$ws.Range('E47').Value = '  -11.44%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0855'
$ws.Range('E48').Value = '  -2.74%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0213'
$ws.Range('E49').Value = '  -5.30%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '16.21'
$ws.Range('E50').Value = '  -13.79%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.10'
$ws.Range('E51').Value = '  -12.30%  '
